$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date header in B1. Prefix with an apostrophe so Excel keeps it
# as literal text instead of auto-converting the "dd/mm/yyyy"-looking string
# into a date serial number/value.
$ws.Range("B1").Value = "'11/03/2023"

# Update the numeric values in B2:B17
$ws.Range("B2").Value = 66
$ws.Range("B3").Value = 136
$ws.Range("B4").Value = 118
$ws.Range("B5").Value = 110
$ws.Range("B6").Value = 101
$ws.Range("B7").Value = 87
$ws.Range("B8").Value = 84
$ws.Range("B9").Value = 106
$ws.Range("B10").Value = 88
$ws.Range("B11").Value = 91
$ws.Range("B12").Value = 77
$ws.Range("B13").Value = 68
$ws.Range("B14").Value = 53
$ws.Range("B15").Value = 22
$ws.Range("B16").Value = 20
$ws.Range("B17").Value = 16
